# BOM update: high-voltage circuitry parts (SHV cable ends soldered to SMD
# pads, potted + zip-tie strain relief) plus SMA ports. Fills in rows 3-8
# of the BOM with component label / type / value / vendor / product number
# / "here" datasheet-or-vendor-link hyperlinks, and marks the "todo" status
# column (H) for each of these new rows (including row 2).
#
# Cell values are written in the same order the original author typed them
# (inferred from shared-string first-use order) so the workbook's shared
# string table comes out in the same sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: SMD HV resistor, 2512 (3M, Mouser 660-HV733ATTE3004F) ----------
$ws.Range("M3").Value = "here"
$ws.Range("L3").Value = "660-HV733ATTE3004F"
$ws.Range("K3").Value = "Mouser"
$ws.Range("J3").Value = "3M"
$ws.Range("N3").Value = "For example, 2512 package type HV resistors"

# --- Row 4: SMD HV cap, 1812 (100 pF, Mouser 581-1812HA101k) ---------------
$ws.Range("J4").Value = "100 pF"
$ws.Range("I3").Value = "SMD HV resistor, 2512"
$ws.Range("L4").Value = "581-1812HA101k"

# --- Row 5: SMD HV cap, 2220 (10 nF, Mouser 80-C2220C103KHR) ---------------
$ws.Range("J5").Value = "10 nF"
$ws.Range("I5").Value = "SMD HV cap, 2220 "
$ws.Range("I4").Value = "SMD HV cap, 1812"

# --- "todo" markers on the status column (H) for the new rows -------------
$ws.Range("H2").Value = "todo"
$ws.Range("H3").Value = "todo"
$ws.Range("H4").Value = "todo"
$ws.Range("H5").Value = "todo"
$ws.Range("H6").Value = "todo"
$ws.Range("H7").Value = "todo"
$ws.Range("H8").Value = "todo"

# --- Row 5 product number, Row 6: 200V SMD Zener (Mouser 78-BZD27C200P) ---
$ws.Range("L5").Value = "80-C2220C103KHR"
$ws.Range("I6").Value = "200V SMD Zener"
$ws.Range("N6").Value = "D_SOD123"
$ws.Range("J6").Value = "200V Zener"
$ws.Range("L6").Value = "78-BZD27C200P-HE3-08"

# --- Row 7: SMD gnd cap, 1812 (0.1 uF, Mouser) -----------------------------
$ws.Range("J7").Value = "0.1 uF"
$ws.Range("I7").Value = "SMD gnd cap, 1812"

# --- Row 8: Angled SMD SMA Jack (Digikey 142-0711-271) ---------------------
$ws.Range("I8").Value = "Angled SMD SMA Jack"
$ws.Range("K8").Value = "Digikey"
$ws.Range("L8").Value = "142-0711-271"

# --- vendor cells that reuse the "Mouser" shared string --------------------
$ws.Range("K4").Value = "Mouser"
$ws.Range("K5").Value = "Mouser"
$ws.Range("K6").Value = "Mouser"
$ws.Range("K7").Value = "Mouser"

# --- "here" hyperlink cells reusing the shared string, plus live links ----
$ws.Range("M4").Value = "here"
$ws.Range("M5").Value = "here"
$ws.Range("M6").Value = "here"
$ws.Range("M7").Value = "here"
$ws.Range("M8").Value = "here"

$ws.Hyperlinks.Add($ws.Range("M3"), "https://www.mouser.com/ProductDetail/660-HV733ATTE3004F") | Out-Null
$ws.Hyperlinks.Add($ws.Range("M4"), "https://www.mouser.com/ProductDetail/581-1812HA101k") | Out-Null
$ws.Hyperlinks.Add($ws.Range("M5"), "https://www.mouser.com/ProductDetail/80-C2220C103KHR") | Out-Null
$ws.Hyperlinks.Add($ws.Range("M6"), "https://www.mouser.com/ProductDetail/78-BZD27C200P-HE3-08") | Out-Null
$ws.Hyperlinks.Add($ws.Range("M7"), "https://www.mouser.com/c/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("M8"), "https://www.digikey.com/en/products/detail/142-0711-271") | Out-Null

# Re-apply the Hyperlink cell style: Hyperlinks.Add stamps the "Hyperlink"
# built-in style on these cells, which is what they already carried (M2:M30
# were pre-styled for links before any links existed) -- keep them on that
# same style after linking.
$ws.Range("M3").Style = "Hyperlink"
$ws.Range("M4").Style = "Hyperlink"
$ws.Range("M5").Style = "Hyperlink"
$ws.Range("M6").Style = "Hyperlink"
$ws.Range("M7").Style = "Hyperlink"
$ws.Range("M8").Style = "Hyperlink"

# Final cursor position left on J8, matching the saved selection.
$ws.Range("J8").Select() | Out-Null

# Best-effort: restore the window geometry recorded in the workbook view.
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 6780
    $win.Top = 23540
    $win.Width = 28800
    $win.Height = 16260
} catch {
}
